$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Copy formatting from the last existing data row (78) down into the new
# row (79) so the new cells inherit the date/number styling used by the
# rest of the table, then overwrite with the actual values.
$ws.Range("A78:F78").Copy()
$ws.Range("A79:F79").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A79").Value = 43981
$ws.Range("B79").Value = 561
$ws.Range("C79").Value = 208
$ws.Range("D79").Value = 471
$ws.Range("E79").Value = 13
$ws.Range("F79").Value = 40

# Grow the table so it covers the newly appended row
$lo = $ws.ListObjects.Item("Condicion_Pacientes")
$lo.Resize($ws.Range("A1:F79"))

# Match the resulting selection recorded in the saved workbook
$ws.Range("D79").Select()
